$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "BRG_01_03",
    "BRG_16_01",
    "BRG_05_03",
    "BRG_05_14",
    "BRG_16_08",
    "BRG_05_09",
    "ECO_09_05",
    "BRG_16_04B",
    "BRG_01_01",
    "BRG_01_04",
    "BRG_05_02",
    "BRG_01_07",
    "BRG_05_04",
    "BRG_16_09",
    "BRG_01_05",
    "ECO_09_01",
    "BRG_16_02",
    "ECO_09_02",
    "BRG_05_05",
    "BRG_01_08",
    "BRG_16_05",
    "BRG_16_03",
    "BRG_05_01",
    "BRG_05_13",
    "BRG_01_09",
    "BRG_13_02",
    "BRG_01_02",
    "BRG_01_06"
)

# New "hole_id" header in A1, matching the bold/bordered header style already
# used by the other header cells (B1:M1).
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Replace the numeric index column (0..27) with the hole_id string for each row.
for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
